$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2357
$ws.Range("I40").Value = 1945
$ws.Range("J40").Value = 2975
$ws.Range("K40").Value = 1945
$ws.Range("L40").Value = 2975
$ws.Range("M40").Value = -1770
$ws.Range("N40").Value = -3325
$ws.Range("H49").Value = 995.6875
$ws.Range("I49").Value = 1042.125
$ws.Range("J49").Value = 949.25
$ws.Range("K49").Value = 3126.375
$ws.Range("L49").Value = 2847.75
$ws.Range("M49").Value = -2990.375
$ws.Range("N49").Value = -3119.75
$ws.Range("H59").Value = 904.875
$ws.Range("J59").Value = 904.875
$ws.Range("L59").Value = 2714.625
$ws.Range("N59").Value = -3828.625
$ws.Range("H61").Value = 319.33334
$ws.Range("I61").Value = 199.8
$ws.Range("K61").Value = 599.4000000000001
$ws.Range("M61").Value = -427.4000000000001
$ws.Range("H105").Value = 40000
$ws.Range("J105").Value = 40000
$ws.Range("L105").Value = 40000
$ws.Range("N105").Value = -46988
$ws.Range("H113").Value = 2551.7917
$ws.Range("I113").Value = 2014.2858
$ws.Range("J113").Value = 2773.1177
$ws.Range("K113").Value = 2014.2858
$ws.Range("L113").Value = 2773.1177
$ws.Range("M113").Value = 1239.7142
$ws.Range("N113").Value = -9281.117699999999
$ws.Range("H125").Value = 349.23077
$ws.Range("I125").Value = 271.42856
$ws.Range("J125").Value = 440
$ws.Range("K125").Value = 2442.85704
$ws.Range("L125").Value = 3960
$ws.Range("M125").Value = 17.14296000000013
$ws.Range("N125").Value = -8880
$ws.Range("H141").Value = 3913.84
$ws.Range("I141").Value = 1804.3513
$ws.Range("J141").Value = 9917.77
$ws.Range("K141").Value = 5413.0539
$ws.Range("L141").Value = 29753.31
$ws.Range("M141").Value = -233.0538999999999
$ws.Range("N141").Value = -40113.31

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 670.0833
$ws.Range("I2").Value = 658.8
$ws.Range("J2").Value = 726.5
$ws.Range("K2").Value = 658.8
$ws.Range("L2").Value = 726.5
$ws.Range("M2").Value = -545.8
$ws.Range("N2").Value = -952.5
$ws.Range("H61").Value = 925.4318
$ws.Range("I61").Value = 815.71875
$ws.Range("J61").Value = 1218
$ws.Range("K61").Value = 815.71875
$ws.Range("L61").Value = 1218
$ws.Range("M61").Value = -603.71875
$ws.Range("N61").Value = -1642
$ws.Range("H74").Value = 5103938.5
$ws.Range("I74").Value = 5953899.5
$ws.Range("J74").Value = 4173.857
$ws.Range("K74").Value = 5953899.5
$ws.Range("L74").Value = 4173.857
$ws.Range("M74").Value = -5953025.5
$ws.Range("N74").Value = -5921.857
$ws.Range("H77").Value = 5103938.5
$ws.Range("I77").Value = 5953899.5
$ws.Range("J77").Value = 4173.857
$ws.Range("K77").Value = 29769497.5
$ws.Range("L77").Value = 20869.285
$ws.Range("M77").Value = -29765129.5
$ws.Range("N77").Value = -29605.285
$ws.Range("H116").Value = 670.0833
$ws.Range("I116").Value = 658.8
$ws.Range("J116").Value = 726.5
$ws.Range("K116").Value = 658.8
$ws.Range("L116").Value = 726.5
$ws.Range("M116").Value = 1635.2
$ws.Range("N116").Value = -5314.5
$ws.Range("H136").Value = 925.4318
$ws.Range("I136").Value = 815.71875
$ws.Range("J136").Value = 1218
$ws.Range("K136").Value = 2447.15625
$ws.Range("L136").Value = 3654
$ws.Range("M136").Value = 102.84375
$ws.Range("N136").Value = -8754

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 670.0833
$ws.Range("I3").Value = 658.8
$ws.Range("J3").Value = 726.5
$ws.Range("K3").Value = 658.8
$ws.Range("L3").Value = 726.5
$ws.Range("M3").Value = -544.8
$ws.Range("N3").Value = -954.5
$ws.Range("H134").Value = 1134.18
$ws.Range("I134").Value = 1050.225
$ws.Range("J134").Value = 1470
$ws.Range("K134").Value = 3150.675
$ws.Range("L134").Value = 4410
$ws.Range("M134").Value = -615.6749999999997
$ws.Range("N134").Value = -9480

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 18000
$ws.Range("J43").Value = 18000
$ws.Range("L43").Value = 18000
$ws.Range("N43").Value = -18368
$ws.Range("H96").Value = 18809.666
$ws.Range("J96").Value = 18809.666
$ws.Range("L96").Value = 18809.666
$ws.Range("N96").Value = -24301.666
$ws.Range("H99").Value = 3577896.2
$ws.Range("I99").Value = 4470745.5
$ws.Range("J99").Value = 6500
$ws.Range("K99").Value = 4470745.5
$ws.Range("L99").Value = 6500
$ws.Range("M99").Value = -4469247.5
$ws.Range("N99").Value = -9496
$ws.Range("H101").Value = 18000
$ws.Range("J101").Value = 18000
$ws.Range("L101").Value = 18000
$ws.Range("N101").Value = -24490
$ws.Range("H106").Value = 30111
$ws.Range("J106").Value = 30111
$ws.Range("L106").Value = 30111
$ws.Range("N106").Value = -32635
$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680
$ws.Range("H126").Value = 3577896.2
$ws.Range("I126").Value = 4470745.5
$ws.Range("J126").Value = 6500
$ws.Range("K126").Value = 13412236.5
$ws.Range("L126").Value = 19500
$ws.Range("M126").Value = -13409766.5
$ws.Range("N126").Value = -24440
$ws.Range("H134").Value = 2647.5454
$ws.Range("I134").Value = 3032.261
$ws.Range("J134").Value = 1762.7
$ws.Range("K134").Value = 9096.782999999999
$ws.Range("L134").Value = 5288.1
$ws.Range("M134").Value = -6561.782999999999
$ws.Range("N134").Value = -10358.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1091.3684
$ws.Range("J131").Value = 1195.3125
$ws.Range("L131").Value = 3585.9375
$ws.Range("N131").Value = -13665.9375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 39000
$ws.Range("J101").Value = 39000
$ws.Range("L101").Value = 39000
$ws.Range("N101").Value = -45490
$ws.Range("H107").Value = 770.5417
$ws.Range("I107").Value = 708
$ws.Range("J107").Value = 833.0833
$ws.Range("K107").Value = 708
$ws.Range("L107").Value = 833.0833
$ws.Range("M107").Value = 1212
$ws.Range("N107").Value = -4673.0833
$ws.Range("H126").Value = 2635.9167
$ws.Range("I126").Value = 1769.8667
$ws.Range("K126").Value = 5309.6001
$ws.Range("M126").Value = -2839.6001
$ws.Range("H132").Value = 1533.5
$ws.Range("I132").Value = 953
$ws.Range("J132").Value = 1920.5
$ws.Range("K132").Value = 2859
$ws.Range("L132").Value = 5761.5
$ws.Range("M132").Value = -329
$ws.Range("N132").Value = -10821.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 306.05884
$ws.Range("I55").Value = 274.3
$ws.Range("J55").Value = 351.42856
$ws.Range("K55").Value = 274.3
$ws.Range("L55").Value = 351.42856
$ws.Range("M55").Value = -101.3
$ws.Range("N55").Value = -697.4285600000001
$ws.Range("H132").Value = 7915425.5
$ws.Range("I132").Value = 17863954
$ws.Range("J132").Value = 1823.0227
$ws.Range("K132").Value = 53591862
$ws.Range("L132").Value = 5469.0681
$ws.Range("M132").Value = -53589332
$ws.Range("N132").Value = -10529.0681
$ws.Range("H136").Value = 2863.1692
$ws.Range("I136").Value = 3324.3618
$ws.Range("J136").Value = 1658.9445
$ws.Range("K136").Value = 9973.0854
$ws.Range("L136").Value = 4976.833500000001
$ws.Range("M136").Value = -7423.0854
$ws.Range("N136").Value = -10076.8335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1248.5834
$ws.Range("I132").Value = 651.73334
$ws.Range("J132").Value = 4232.8335
$ws.Range("K132").Value = 1955.20002
$ws.Range("L132").Value = 12698.5005
$ws.Range("M132").Value = 574.79998
$ws.Range("N132").Value = -17758.5005
$ws.Range("H136").Value = 872.3570999999999
$ws.Range("I136").Value = 573.2712
$ws.Range("J136").Value = 1578.2
$ws.Range("K136").Value = 1719.8136
$ws.Range("L136").Value = 4734.6
$ws.Range("M136").Value = 830.1864
$ws.Range("N136").Value = -9834.6
